$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Insert two new data rows (18 and 19) right below the current last
#    data row (17). This pushes the trailing signature block (old rows
#    22-23) down to rows 24-25 automatically.
# ---------------------------------------------------------------------
$ws.Rows("18:19").Insert()

# Preserve the heavier "closing" border that used to belong to row 17
# (the old last row) by copying its formatting down onto the new last
# row (19) before row 17 itself gets reformatted. Use the bounded table
# range (B:J) so the untouched columns outside the table are not
# disturbed.
$ws.Range("B17:J17").Copy()
$ws.Range("B19:J19").PasteSpecial(-4122)

# Row 17 becomes a regular middle row now, so it should pick up the
# same formatting as row 16 (the other middle row). Stamp that format
# onto rows 17 and 18 as well.
$ws.Range("B16:J16").Copy()
$ws.Range("B17:J18").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 2. New shared string used for the second period column value.
# ---------------------------------------------------------------------
$periodo2 = "2508"

# ---------------------------------------------------------------------
# 3. Update the existing data rows' amounts for the new period totals.
# ---------------------------------------------------------------------
$ws.Range("G16").Value = 1500000

$ws.Range("B17").Value = "CC"
$ws.Range("C17").Value = "18857167"
$ws.Range("D17").Value = "YONIS DEL CRISTO LARIO ALDANA"
$ws.Range("E17").Value = "2505"
$ws.Range("F17").Value = 56940
$ws.Range("G17").Value = 1423500

# ---------------------------------------------------------------------
# 4. Fill in the two brand-new rows for period 2508.
# ---------------------------------------------------------------------
$ws.Range("B18").Value = "CC"
$ws.Range("C18").Value = "92670739"
$ws.Range("D18").Value = "MIGUEL BELTRAN SALGADO"
$ws.Range("E18").Value = $periodo2
$ws.Range("F18").Value = 60000
$ws.Range("G18").Value = 1500000

$ws.Range("B19").Value = "CC"
$ws.Range("C19").Value = "18857167"
$ws.Range("D19").Value = "YONIS DEL CRISTO LARIO ALDANA"
$ws.Range("E19").Value = $periodo2
$ws.Range("F19").Value = 56940
$ws.Range("G19").Value = 1423500

# ---------------------------------------------------------------------
# 5. Refresh the summary block above the table: total overdue value and
#    the number of periods now covered.
# ---------------------------------------------------------------------
$ws.Range("E11").Value = 233880
$ws.Range("F13").Value = 2
